$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,5).Value = 3
$ws.Cells.Item(2,6).Value = 1
$ws.Cells.Item(2,7).Value = 19.23438166666667
$ws.Cells.Item(2,8).Value = 57.70314500000001
$ws.Cells.Item(2,9).Value = 0.1061888747949397
$ws.Cells.Item(2,10).Value = 0.1061888747949397
$ws.Cells.Item(2,11).Value = 3
$ws.Cells.Item(2,12).Value = 1
$ws.Cells.Item(2,13).Value = 9.519603666666667
$ws.Cells.Item(2,14).Value = 28.558811
$ws.Cells.Item(2,15).Value = 0.1143733649354437
$ws.Cells.Item(2,16).Value = 0.1143733649354437
$ws.Cells.Item(2,17).Value = 183.1036902400661
$ws.Cells.Item(2,18).Value = 1647.933212160595
$ws.Cells.Item(2,19).Value = 0.01214517892900578
$ws.Cells.Item(2,20).Value = 0.01214517892900577

$ws.Cells.Item(3,5).Value = 3
$ws.Cells.Item(3,6).Value = 1
$ws.Cells.Item(3,7).Value = 19.23438166666667
$ws.Cells.Item(3,8).Value = 57.70314500000001
$ws.Cells.Item(3,9).Value = 0.1061888747949397
$ws.Cells.Item(3,10).Value = 0.1061888747949397
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 3.397950666666667
$ws.Cells.Item(3,14).Value = 10.193852
$ws.Cells.Item(3,15).Value = 0.04082470922525111
$ws.Cells.Item(3,16).Value = 0.04082470922525111
$ws.Cells.Item(3,17).Value = 65.3574800071711
$ws.Cells.Item(3,18).Value = 588.21732006454
$ws.Cells.Item(3,19).Value = 0.004335129936460009
$ws.Cells.Item(3,20).Value = 0.004335129936460009

$ws.Cells.Item(4,5).Value = 3
$ws.Cells.Item(4,6).Value = 1
$ws.Cells.Item(4,7).Value = 19.23438166666667
$ws.Cells.Item(4,8).Value = 57.70314500000001
$ws.Cells.Item(4,9).Value = 0.1061888747949397
$ws.Cells.Item(4,10).Value = 0.1061888747949397
$ws.Cells.Item(4,11).Value = 3
$ws.Cells.Item(4,12).Value = 1
$ws.Cells.Item(4,13).Value = 35.04087533333333
$ws.Cells.Item(4,14).Value = 105.122626
$ws.Cells.Item(4,15).Value = 0.4209989157626403
$ws.Cells.Item(4,16).Value = 0.4209989157626403
$ws.Cells.Item(4,17).Value = 673.9895700954189
$ws.Cells.Item(4,18).Value = 6065.906130858771
$ws.Cells.Item(4,19).Value = 0.04470540115472436
$ws.Cells.Item(4,20).Value = 0.04470540115472436

$ws.Cells.Item(5,5).Value = 3
$ws.Cells.Item(5,6).Value = 1
$ws.Cells.Item(5,7).Value = 19.23438166666667
$ws.Cells.Item(5,8).Value = 57.70314500000001
$ws.Cells.Item(5,9).Value = 0.1061888747949397
$ws.Cells.Item(5,10).Value = 0.1061888747949397
$ws.Cells.Item(5,11).Value = 3
$ws.Cells.Item(5,12).Value = 1
$ws.Cells.Item(5,13).Value = 34.43746266666667
$ws.Cells.Item(5,14).Value = 103.312388
$ws.Cells.Item(5,15).Value = 0.4137492087845029
$ws.Cells.Item(5,16).Value = 0.4137492087845028
$ws.Cells.Item(5,17).Value = 662.3833005622512
$ws.Cells.Item(5,18).Value = 5961.44970506026
$ws.Cells.Item(5,19).Value = 0.04393556292812293
$ws.Cells.Item(5,20).Value = 0.04393556292812292

$ws.Cells.Item(6,5).Value = 3
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = 19.23438166666667
$ws.Cells.Item(6,8).Value = 57.70314500000001
$ws.Cells.Item(6,9).Value = 0.1061888747949397
$ws.Cells.Item(6,10).Value = 0.1061888747949397
$ws.Cells.Item(6,11).Value = 3
$ws.Cells.Item(6,12).Value = 1
$ws.Cells.Item(6,13).Value = 0.836805
$ws.Cells.Item(6,14).Value = 2.510415
$ws.Cells.Item(6,15).Value = 0.01005380129216206
$ws.Cells.Item(6,16).Value = 0.01005380129216206
$ws.Cells.Item(6,17).Value = 16.095426750575
$ws.Cells.Item(6,18).Value = 144.858840755175
$ws.Cells.Item(6,19).Value = 0.001067601846626599
$ws.Cells.Item(6,20).Value = 0.001067601846626599

$ws.Cells.Item(7,5).Value = 3
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(7,7).Value = 101.9328183333333
$ws.Cells.Item(7,8).Value = 305.798455
$ws.Cells.Item(7,9).Value = 0.5627491161960234
$ws.Cells.Item(7,10).Value = 0.5627491161960234
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 9.519603666666667
$ws.Cells.Item(7,14).Value = 28.558811
$ws.Cells.Item(7,15).Value = 0.1143733649354437
$ws.Cells.Item(7,16).Value = 0.1143733649354437
$ws.Cells.Item(7,17).Value = 970.3600311596672
$ws.Cells.Item(7,18).Value = 8733.240280437005
$ws.Cells.Item(7,19).Value = 0.06436351003378621
$ws.Cells.Item(7,20).Value = 0.06436351003378621

$ws.Cells.Item(8,5).Value = 3
$ws.Cells.Item(8,6).Value = 1
$ws.Cells.Item(8,7).Value = 101.9328183333333
$ws.Cells.Item(8,8).Value = 305.798455
$ws.Cells.Item(8,9).Value = 0.5627491161960234
$ws.Cells.Item(8,10).Value = 0.5627491161960234
$ws.Cells.Item(8,11).Value = 3
$ws.Cells.Item(8,12).Value = 1
$ws.Cells.Item(8,13).Value = 3.397950666666667
$ws.Cells.Item(8,14).Value = 10.193852
$ws.Cells.Item(8,15).Value = 0.04082470922525111
$ws.Cells.Item(8,16).Value = 0.04082470922525111
$ws.Cells.Item(8,17).Value = 346.3626880109622
$ws.Cells.Item(8,18).Value = 3117.26419209866
$ws.Cells.Item(8,19).Value = 0.02297406903546971
$ws.Cells.Item(8,20).Value = 0.02297406903546971

$ws.Cells.Item(9,5).Value = 3
$ws.Cells.Item(9,6).Value = 1
$ws.Cells.Item(9,7).Value = 101.9328183333333
$ws.Cells.Item(9,8).Value = 305.798455
$ws.Cells.Item(9,9).Value = 0.5627491161960234
$ws.Cells.Item(9,10).Value = 0.5627491161960234
$ws.Cells.Item(9,11).Value = 3
$ws.Cells.Item(9,12).Value = 1
$ws.Cells.Item(9,13).Value = 35.04087533333333
$ws.Cells.Item(9,14).Value = 105.122626
$ws.Cells.Item(9,15).Value = 0.4209989157626403
$ws.Cells.Item(9,16).Value = 0.4209989157626403
$ws.Cells.Item(9,17).Value = 3571.815179593647
$ws.Cells.Item(9,18).Value = 32146.33661634283
$ws.Cells.Item(9,19).Value = 0.2369167677649099
$ws.Cells.Item(9,20).Value = 0.2369167677649099

$ws.Cells.Item(10,5).Value = 3
$ws.Cells.Item(10,6).Value = 1
$ws.Cells.Item(10,7).Value = 101.9328183333333
$ws.Cells.Item(10,8).Value = 305.798455
$ws.Cells.Item(10,9).Value = 0.5627491161960234
$ws.Cells.Item(10,10).Value = 0.5627491161960234
$ws.Cells.Item(10,11).Value = 3
$ws.Cells.Item(10,12).Value = 1
$ws.Cells.Item(10,13).Value = 34.43746266666667
$ws.Cells.Item(10,14).Value = 103.312388
$ws.Cells.Item(10,15).Value = 0.4137492087845029
$ws.Cells.Item(10,16).Value = 0.4137492087845028
$ws.Cells.Item(10,17).Value = 3510.307625862282
$ws.Cells.Item(10,18).Value = 31592.76863276054
$ws.Cells.Item(10,19).Value = 0.232837001570283
$ws.Cells.Item(10,20).Value = 0.2328370015702829

$ws.Cells.Item(11,5).Value = 3
$ws.Cells.Item(11,6).Value = 1
$ws.Cells.Item(11,7).Value = 101.9328183333333
$ws.Cells.Item(11,8).Value = 305.798455
$ws.Cells.Item(11,9).Value = 0.5627491161960234
$ws.Cells.Item(11,10).Value = 0.5627491161960234
$ws.Cells.Item(11,11).Value = 3
$ws.Cells.Item(11,12).Value = 1
$ws.Cells.Item(11,13).Value = 0.836805
$ws.Cells.Item(11,14).Value = 2.510415
$ws.Cells.Item(11,15).Value = 0.01005380129216206
$ws.Cells.Item(11,16).Value = 0.01005380129216206
$ws.Cells.Item(11,17).Value = 85.297892045425
$ws.Cells.Item(11,18).Value = 767.681028408825
$ws.Cells.Item(11,19).Value = 0.005657767791574635
$ws.Cells.Item(11,20).Value = 0.005657767791574635

$ws.Cells.Item(12,5).Value = 3
$ws.Cells.Item(12,6).Value = 1
$ws.Cells.Item(12,7).Value = 33.990832
$ws.Cells.Item(12,8).Value = 101.972496
$ws.Cells.Item(12,9).Value = 0.1876560560134372
$ws.Cells.Item(12,10).Value = 0.1876560560134371
$ws.Cells.Item(12,11).Value = 3
$ws.Cells.Item(12,12).Value = 1
$ws.Cells.Item(12,13).Value = 9.519603666666667
$ws.Cells.Item(12,14).Value = 28.558811
$ws.Cells.Item(12,15).Value = 0.1143733649354437
$ws.Cells.Item(12,16).Value = 0.1143733649354437
$ws.Cells.Item(12,17).Value = 323.5792489402507
$ws.Cells.Item(12,18).Value = 2912.213240462256
$ws.Cells.Item(12,19).Value = 0.02146285457677092
$ws.Cells.Item(12,20).Value = 0.02146285457677091

$ws.Cells.Item(13,5).Value = 3
$ws.Cells.Item(13,6).Value = 1
$ws.Cells.Item(13,7).Value = 33.990832
$ws.Cells.Item(13,8).Value = 101.972496
$ws.Cells.Item(13,9).Value = 0.1876560560134372
$ws.Cells.Item(13,10).Value = 0.1876560560134371
$ws.Cells.Item(13,11).Value = 3
$ws.Cells.Item(13,12).Value = 1
$ws.Cells.Item(13,13).Value = 3.397950666666667
$ws.Cells.Item(13,14).Value = 10.193852
$ws.Cells.Item(13,15).Value = 0.04082470922525111
$ws.Cells.Item(13,16).Value = 0.04082470922525111
$ws.Cells.Item(13,17).Value = 115.4991702549547
$ws.Cells.Item(13,18).Value = 1039.492532294592
$ws.Cells.Item(13,19).Value = 0.007661003921106008
$ws.Cells.Item(13,20).Value = 0.007661003921106006

$ws.Cells.Item(14,5).Value = 3
$ws.Cells.Item(14,6).Value = 1
$ws.Cells.Item(14,7).Value = 33.990832
$ws.Cells.Item(14,8).Value = 101.972496
$ws.Cells.Item(14,9).Value = 0.1876560560134372
$ws.Cells.Item(14,10).Value = 0.1876560560134371
$ws.Cells.Item(14,11).Value = 3
$ws.Cells.Item(14,12).Value = 1
$ws.Cells.Item(14,13).Value = 35.04087533333333
$ws.Cells.Item(14,14).Value = 105.122626
$ws.Cells.Item(14,15).Value = 0.4209989157626403
$ws.Cells.Item(14,16).Value = 0.4209989157626403
$ws.Cells.Item(14,17).Value = 1191.068506588277
$ws.Cells.Item(14,18).Value = 10719.6165592945
$ws.Cells.Item(14,19).Value = 0.07900299611795034
$ws.Cells.Item(14,20).Value = 0.07900299611795032

$ws.Cells.Item(15,5).Value = 3
$ws.Cells.Item(15,6).Value = 1
$ws.Cells.Item(15,7).Value = 33.990832
$ws.Cells.Item(15,8).Value = 101.972496
$ws.Cells.Item(15,9).Value = 0.1876560560134372
$ws.Cells.Item(15,10).Value = 0.1876560560134371
$ws.Cells.Item(15,11).Value = 3
$ws.Cells.Item(15,12).Value = 1
$ws.Cells.Item(15,13).Value = 34.43746266666667
$ws.Cells.Item(15,14).Value = 103.312388
$ws.Cells.Item(15,15).Value = 0.4137492087845029
$ws.Cells.Item(15,16).Value = 0.4137492087845028
$ws.Cells.Item(15,17).Value = 1170.558008008939
$ws.Cells.Item(15,18).Value = 10535.02207208045
$ws.Cells.Item(15,19).Value = 0.07764254469917999
$ws.Cells.Item(15,20).Value = 0.07764254469917996

$ws.Cells.Item(16,5).Value = 3
$ws.Cells.Item(16,6).Value = 1
$ws.Cells.Item(16,7).Value = 33.990832
$ws.Cells.Item(16,8).Value = 101.972496
$ws.Cells.Item(16,9).Value = 0.1876560560134372
$ws.Cells.Item(16,10).Value = 0.1876560560134371
$ws.Cells.Item(16,11).Value = 3
$ws.Cells.Item(16,12).Value = 1
$ws.Cells.Item(16,13).Value = 0.836805
$ws.Cells.Item(16,14).Value = 2.510415
$ws.Cells.Item(16,15).Value = 0.01005380129216206
$ws.Cells.Item(16,16).Value = 0.01005380129216206
$ws.Cells.Item(16,17).Value = 28.44369817176
$ws.Cells.Item(16,18).Value = 255.99328354584
$ws.Cells.Item(16,19).Value = 0.00188665669842993
$ws.Cells.Item(16,20).Value = 0.001886656698429929

$ws.Cells.Item(17,5).Value = 3
$ws.Cells.Item(17,6).Value = 1
$ws.Cells.Item(17,7).Value = 24.872162
$ws.Cells.Item(17,8).Value = 74.61648600000001
$ws.Cells.Item(17,9).Value = 0.1373138446698593
$ws.Cells.Item(17,10).Value = 0.1373138446698593
$ws.Cells.Item(17,11).Value = 3
$ws.Cells.Item(17,12).Value = 1
$ws.Cells.Item(17,13).Value = 9.519603666666667
$ws.Cells.Item(17,14).Value = 28.558811
$ws.Cells.Item(17,15).Value = 0.1143733649354437
$ws.Cells.Item(17,16).Value = 0.1143733649354437
$ws.Cells.Item(17,17).Value = 236.7731245731274
$ws.Cells.Item(17,18).Value = 2130.958121158146
$ws.Cells.Item(17,19).Value = 0.01570504646711465
$ws.Cells.Item(17,20).Value = 0.01570504646711465

$ws.Cells.Item(18,5).Value = 3
$ws.Cells.Item(18,6).Value = 1
$ws.Cells.Item(18,7).Value = 24.872162
$ws.Cells.Item(18,8).Value = 74.61648600000001
$ws.Cells.Item(18,9).Value = 0.1373138446698593
$ws.Cells.Item(18,10).Value = 0.1373138446698593
$ws.Cells.Item(18,11).Value = 3
$ws.Cells.Item(18,12).Value = 1
$ws.Cells.Item(18,13).Value = 3.397950666666667
$ws.Cells.Item(18,14).Value = 10.193852
$ws.Cells.Item(18,15).Value = 0.04082470922525111
$ws.Cells.Item(18,16).Value = 0.04082470922525111
$ws.Cells.Item(18,17).Value = 84.51437944934133
$ws.Cells.Item(18,18).Value = 760.6294150440721
$ws.Cells.Item(18,19).Value = 0.005605797781248304
$ws.Cells.Item(18,20).Value = 0.005605797781248303

$ws.Cells.Item(19,5).Value = 3
$ws.Cells.Item(19,6).Value = 1
$ws.Cells.Item(19,7).Value = 24.872162
$ws.Cells.Item(19,8).Value = 74.61648600000001
$ws.Cells.Item(19,9).Value = 0.1373138446698593
$ws.Cells.Item(19,10).Value = 0.1373138446698593
$ws.Cells.Item(19,11).Value = 3
$ws.Cells.Item(19,12).Value = 1
$ws.Cells.Item(19,13).Value = 35.04087533333333
$ws.Cells.Item(19,14).Value = 105.122626
$ws.Cells.Item(19,15).Value = 0.4209989157626403
$ws.Cells.Item(19,16).Value = 0.4209989157626403
$ws.Cells.Item(19,17).Value = 871.5423279124708
$ws.Cells.Item(19,18).Value = 7843.880951212237
$ws.Cells.Item(19,19).Value = 0.05780897972521037
$ws.Cells.Item(19,20).Value = 0.05780897972521036

$ws.Cells.Item(20,5).Value = 3
$ws.Cells.Item(20,6).Value = 1
$ws.Cells.Item(20,7).Value = 24.872162
$ws.Cells.Item(20,8).Value = 74.61648600000001
$ws.Cells.Item(20,9).Value = 0.1373138446698593
$ws.Cells.Item(20,10).Value = 0.1373138446698593
$ws.Cells.Item(20,11).Value = 3
$ws.Cells.Item(20,12).Value = 1
$ws.Cells.Item(20,13).Value = 34.43746266666667
$ws.Cells.Item(20,14).Value = 103.312388
$ws.Cells.Item(20,15).Value = 0.4137492087845029
$ws.Cells.Item(20,16).Value = 0.4137492087845028
$ws.Cells.Item(20,17).Value = 856.5341503142855
$ws.Cells.Item(20,18).Value = 7708.807352828569
$ws.Cells.Item(20,19).Value = 0.05681349458731243
$ws.Cells.Item(20,20).Value = 0.05681349458731241

$ws.Cells.Item(21,5).Value = 3
$ws.Cells.Item(21,6).Value = 1
$ws.Cells.Item(21,7).Value = 24.872162
$ws.Cells.Item(21,8).Value = 74.61648600000001
$ws.Cells.Item(21,9).Value = 0.1373138446698593
$ws.Cells.Item(21,10).Value = 0.1373138446698593
$ws.Cells.Item(21,11).Value = 3
$ws.Cells.Item(21,12).Value = 1
$ws.Cells.Item(21,13).Value = 0.836805
$ws.Cells.Item(21,14).Value = 2.510415
$ws.Cells.Item(21,15).Value = 0.01005380129216206
$ws.Cells.Item(21,16).Value = 0.01005380129216206
$ws.Cells.Item(21,17).Value = 20.81314952241
$ws.Cells.Item(21,18).Value = 187.31834570169
$ws.Cells.Item(21,19).Value = 0.001380526108973571
$ws.Cells.Item(21,20).Value = 0.001380526108973571

$ws.Cells.Item(22,5).Value = 3
$ws.Cells.Item(22,6).Value = 1
$ws.Cells.Item(22,7).Value = 1.103486
$ws.Cells.Item(22,8).Value = 3.310458
$ws.Cells.Item(22,9).Value = 0.006092108325740414
$ws.Cells.Item(22,10).Value = 0.006092108325740414
$ws.Cells.Item(22,11).Value = 3
$ws.Cells.Item(22,12).Value = 1
$ws.Cells.Item(22,13).Value = 9.519603666666667
$ws.Cells.Item(22,14).Value = 28.558811
$ws.Cells.Item(22,15).Value = 0.1143733649354437
$ws.Cells.Item(22,16).Value = 0.1143733649354437
$ws.Cells.Item(22,17).Value = 10.50474937171533
$ws.Cells.Item(22,18).Value = 94.542744345438
$ws.Cells.Item(22,19).Value = 0.0006967749287661634
$ws.Cells.Item(22,20).Value = 0.0006967749287661633

$ws.Cells.Item(23,5).Value = 3
$ws.Cells.Item(23,6).Value = 1
$ws.Cells.Item(23,7).Value = 1.103486
$ws.Cells.Item(23,8).Value = 3.310458
$ws.Cells.Item(23,9).Value = 0.006092108325740414
$ws.Cells.Item(23,10).Value = 0.006092108325740414
$ws.Cells.Item(23,11).Value = 3
$ws.Cells.Item(23,12).Value = 1
$ws.Cells.Item(23,13).Value = 3.397950666666667
$ws.Cells.Item(23,14).Value = 10.193852
$ws.Cells.Item(23,15).Value = 0.04082470922525111
$ws.Cells.Item(23,16).Value = 0.04082470922525111
$ws.Cells.Item(23,17).Value = 3.749590989357333
$ws.Cells.Item(23,18).Value = 33.746318904216
$ws.Cells.Item(23,19).Value = 0.0002487085509670838
$ws.Cells.Item(23,20).Value = 0.0002487085509670837

$ws.Cells.Item(24,5).Value = 3
$ws.Cells.Item(24,6).Value = 1
$ws.Cells.Item(24,7).Value = 1.103486
$ws.Cells.Item(24,8).Value = 3.310458
$ws.Cells.Item(24,9).Value = 0.006092108325740414
$ws.Cells.Item(24,10).Value = 0.006092108325740414
$ws.Cells.Item(24,11).Value = 3
$ws.Cells.Item(24,12).Value = 1
$ws.Cells.Item(24,13).Value = 35.04087533333333
$ws.Cells.Item(24,14).Value = 105.122626
$ws.Cells.Item(24,15).Value = 0.4209989157626403
$ws.Cells.Item(24,16).Value = 0.4209989157626403
$ws.Cells.Item(24,17).Value = 38.66711535807866
$ws.Cells.Item(24,18).Value = 348.004038222708
$ws.Cells.Item(24,19).Value = 0.002564770999845268
$ws.Cells.Item(24,20).Value = 0.002564770999845268

$ws.Cells.Item(25,5).Value = 3
$ws.Cells.Item(25,6).Value = 1
$ws.Cells.Item(25,7).Value = 1.103486
$ws.Cells.Item(25,8).Value = 3.310458
$ws.Cells.Item(25,9).Value = 0.006092108325740414
$ws.Cells.Item(25,10).Value = 0.006092108325740414
$ws.Cells.Item(25,11).Value = 3
$ws.Cells.Item(25,12).Value = 1
$ws.Cells.Item(25,13).Value = 34.43746266666667
$ws.Cells.Item(25,14).Value = 103.312388
$ws.Cells.Item(25,15).Value = 0.4137492087845029
$ws.Cells.Item(25,16).Value = 0.4137492087845028
$ws.Cells.Item(25,17).Value = 38.00125792818933
$ws.Cells.Item(25,18).Value = 342.011321353704
$ws.Cells.Item(25,19).Value = 0.002520604999604579
$ws.Cells.Item(25,20).Value = 0.002520604999604578

$ws.Cells.Item(26,5).Value = 3
$ws.Cells.Item(26,6).Value = 1
$ws.Cells.Item(26,7).Value = 1.103486
$ws.Cells.Item(26,8).Value = 3.310458
$ws.Cells.Item(26,9).Value = 0.006092108325740414
$ws.Cells.Item(26,10).Value = 0.006092108325740414
$ws.Cells.Item(26,11).Value = 3
$ws.Cells.Item(26,12).Value = 1
$ws.Cells.Item(26,13).Value = 0.836805
$ws.Cells.Item(26,14).Value = 2.510415
$ws.Cells.Item(26,15).Value = 0.01005380129216206
$ws.Cells.Item(26,16).Value = 0.01005380129216206
$ws.Cells.Item(26,17).Value = 0.92340260223
$ws.Cells.Item(26,18).Value = 8.31062342007
$ws.Cells.Item(26,19).Value = 0.0000612488465573202
$ws.Cells.Item(26,20).Value = 0.0000612488465573202

